$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2544
$ws1.Range("F7").Value = 1344
$ws1.Range("F8").Value = 1704
$ws1.Range("F11").Value = 2382
$ws1.Range("F12").Value = 492
$ws1.Range("F13").Value = 156
$ws1.Range("F18").Value = 8631
$ws1.Range("F20").Value = 6734
$ws1.Range("F21").Value = 10870
$ws1.Range("F24").Value = 212
$ws1.Range("F25").Value = 297
$ws1.Range("F27").Value = 1295
$ws1.Range("F28").Value = 201
$ws1.Range("F29").Value = 176
$ws1.Range("F30").Value = 2201
$ws1.Range("F31").Value = 92
$ws1.Range("F32").Value = 24
$ws1.Range("F34").Value = 458
$ws1.Range("F35").Value = 441

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 1174

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 34

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 2544
$ws4.Range("F10").Value = 1344
$ws4.Range("F12").Value = 1704
$ws4.Range("F15").Value = 2382
$ws4.Range("F17").Value = 492
$ws4.Range("F18").Value = 156
$ws4.Range("F23").Value = 8631
$ws4.Range("F25").Value = 6734
$ws4.Range("F26").Value = 10870
$ws4.Range("F30").Value = 212
$ws4.Range("F31").Value = 297
$ws4.Range("F36").Value = 201
$ws4.Range("F37").Value = 176
$ws4.Range("F38").Value = 24
$ws4.Range("F46").Value = 441
